$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix C72: it currently holds a stray text value ("`") instead of the
# shared formula used by the rest of that column (C = C_prev * $E$1).
$ws.Range("C72").Formula = "=C71*`$E`$1"

# Force a full recalculation so the cascade of #VALUE! errors in
# D72:E111 and Q1:T1 clears up and resolves to numeric results.
$excel.CalculateFullRebuild()

# Update the active selection to match the saved state of the sheet.
$ws.Range("E87").Select()
